$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: paragraph with "Date designed: ..." (2nd paragraph)
#   - run " Date designed: Nov 28, 2024, 6:10 P" -> "Date designed: Nov 28, 2024, 6:10 P"
#     (drop the leading space, no longer needs xml:space="preserve")
#   - run " " (single space) -> "M" (separate run)
# ------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$rng = $p2.Range
$full = $rng.Text
$pStart = $rng.Start
$idxD = $full.IndexOf("Date designed")
$absD = $pStart + $idxD

$runBStart = $absD - 1
$oldTail = " Date designed: Nov 28, 2024, 6:10 P" + " "
$runsEnd = $runBStart + $oldTail.Length

$oldRange = $d.Range($runBStart, $runsEnd)
$oldRange.Text = ""

$insertPoint1 = $d.Range($runBStart, $runBStart)
$insertPoint1.InsertAfter("Date designed: Nov 28, 2024, 6:10 P")

$afterFirstPos = $runBStart + ("Date designed: Nov 28, 2024, 6:10 P").Length
$insertPoint2 = $d.Range($afterFirstPos, $afterFirstPos)
$insertPoint2.InsertAfter("M")

# ------------------------------------------------------------------
# Change 2: paragraph with "Date conducted: " (3rd paragraph)
#   - spacer run of 30 spaces -> 25 spaces
#   - add new run "Nov 28, 2024, 9:30 PM" right after "Date conducted: "
# ------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$rng3 = $p3.Range
$full3 = $rng3.Text
$pStart3 = $rng3.Start
$pEnd3 = $rng3.End
$idxDC = $full3.IndexOf("Date conducted")
$absDC = $pStart3 + $idxDC

# Remove the trailing "Date conducted: " run so the spacer run becomes
# the last run in the paragraph (needed to edit it without it merging
# into its neighbors).
$trailing = $d.Range($absDC, $pEnd3 - 1)
$trailing.Text = ""

# Shrink the spacer run from 30 to 25 spaces while it sits at the very
# end of the paragraph.
$p3b = $d.Paragraphs.Item(3)
$rngB = $p3b.Range
$pEndB = $rngB.End
$spacesEnd = $pEndB - 1
$spacesStart = $spacesEnd - 30
$spacesRange = $d.Range($spacesStart, $spacesEnd)
$spacesRange.Text = ""
$insertPoint3 = $d.Range($spacesStart, $spacesStart)
$insertPoint3.InsertAfter("                         ")

# Re-add "Date conducted: " (trailing char is a non-breaking space).
$p3c = $d.Paragraphs.Item(3)
$rngC = $p3c.Range
$pEndC = $rngC.End
$insertPoint4 = $d.Range($pEndC - 1, $pEndC - 1)
$insertPoint4.InsertAfter("Date conducted:" + [char]0x00A0)

# Add the new date run.
$p3d = $d.Paragraphs.Item(3)
$rngD = $p3d.Range
$pEndD = $rngD.End
$insertPoint5 = $d.Range($pEndD - 1, $pEndD - 1)
$insertPoint5.InsertAfter("Nov 28, 2024, 9:30 PM")
